$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card16")

$ws.Range("B29:K29").Value = "nan"
$ws.Range("L29").Value = "18/12/2025"
